# Update cryptocurrency price/volume figures per the Tue Apr 23 16:16:53 UTC 2024 data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.608.63"
$ws.Range("E2").Value = "  +0.66%  "
$ws.Range("D3").Value = "3.236.04"
$ws.Range("E3").Value = "  +1.39%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'605.70"
$ws.Range("E5").Value = "  +1.54%  "
$ws.Range("D6").Value = "158.72"
$ws.Range("E6").Value = "  +2.84%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "3.235.64"
$ws.Range("E8").Value = "  +1.39%  "
$ws.Range("E9").Value = "  +1.31%  "
$ws.Range("E10").Value = "  +1.80%  "
$ws.Range("D11").Value = "5.72"
$ws.Range("E11").Value = "  -5.65%  "
$ws.Range("D12").Value = "0.507"
$ws.Range("E12").Value = "  -1.97%  "
$ws.Range("E13").Value = "  +2.89%  "
$ws.Range("D14").Value = "39.17"
$ws.Range("E14").Value = "  -0.28%  "
$ws.Range("D15").Value = "3.763.75"
$ws.Range("E15").Value = "  +1.45%  "
$ws.Range("D16").Value = "66.665.90"
$ws.Range("E16").Value = "  +0.98%  "
$ws.Range("D17").Value = "7.42"
$ws.Range("E17").Value = "  -0.60%  "
$ws.Range("D18").Value = "3.241.25"
$ws.Range("E18").Value = "  +1.77%  "
$ws.Range("E19").Value = "  +1.10%  "
$ws.Range("D20").Value = "510.96"
$ws.Range("E20").Value = "  +0.05%  "
$ws.Range("E21").Value = "  -0.18%  "
$ws.Range("E22").Value = "  -0.57%  "
$ws.Range("D23").Value = "8.06"
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("D24").Value = "'14.80"
$ws.Range("E24").Value = "  -2.40%  "
$ws.Range("D25").Value = "84.92"
$ws.Range("E25").Value = "  -0.09%  "
$ws.Range("E26").Value = "  +0.43%  "
$ws.Range("E27").Value = "  -0.22%  "
$ws.Range("E28").Value = "  -0.87%  "
$ws.Range("D29").Value = "2.39"
$ws.Range("E29").Value = "  +3.62%  "
$ws.Range("D30").Value = "2.99"
$ws.Range("E30").Value = "  +2.91%  "
$ws.Range("D31").Value = "7.05"
$ws.Range("E31").Value = "  -0.01%  "
$ws.Range("D32").Value = "28.27"
$ws.Range("E32").Value = "  +0.08%  "
$ws.Range("E33").Value = "  +0.07%  "
$ws.Range("E34").Value = "  -3.49%  "
$ws.Range("D35").Value = "0.101"
$ws.Range("E35").Value = "  +11.83%  "
$ws.Range("D36").Value = "6.53"
$ws.Range("E36").Value = "  +0.17%  "
$ws.Range("D37").Value = "512.15"
$ws.Range("E37").Value = "  +5.85%  "
$ws.Range("D38").Value = "56.18"
$ws.Range("E38").Value = "  +2.42%  "
$ws.Range("D39").Value = "0.0₃0773"
$ws.Range("E39").Value = "  +17.18%  "
$ws.Range("D40").Value = "0.0422"
$ws.Range("D41").Value = "3.06"
$ws.Range("E41").Value = "  +6.84%  "
$ws.Range("D42").Value = "'0.130"
$ws.Range("E42").Value = "  +6.31%  "
$ws.Range("D43").Value = "8.78"
$ws.Range("E43").Value = "  -1.35%  "
$ws.Range("D44").Value = "'0.300"
$ws.Range("E44").Value = "  -1.02%  "
$ws.Range("D45").Value = "2.49"
$ws.Range("E45").Value = "  +2.59%  "
$ws.Range("D46").Value = "2.879.77"
$ws.Range("E46").Value = "  -1.23%  "
$ws.Range("D47").Value = "28.63"
$ws.Range("E47").Value = "  +0.08%  "
$ws.Range("D48").Value = "2.42"
$ws.Range("E48").Value = "  +4.05%  "
$ws.Range("E49").Value = "  -0.02%  "
$ws.Range("E50").Value = "  -0.21%  "

# Row 51: coin replaced (CoreDAO -> Monero)
$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").Value = "'122.50"
$ws.Range("E51").Value = "  +1.04%  "
